$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44495
$ws.Range("J2").Value = 130

# Row 3
$ws.Range("D3").Value = 44418
$ws.Range("H3").Value = 'Española'
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 16000
$ws.Range("N3").Value = '$/caja 30 unidades'
$ws.Range("P3").Value = 533
$ws.Range("Q3").Value = 30

# Row 4
$ws.Range("D4").Value = 44390
$ws.Range("H4").Value = 'Española'
$ws.Range("J4").Value = 80
$ws.Range("N4").Value = '$/caja 30 unidades'
$ws.Range("P4").Value = 533
$ws.Range("Q4").Value = 30

# Row 5
$ws.Range("D5").Value = 44446
$ws.Range("H5").Value = 'Madrigal'
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("P5").Value = 400
$ws.Range("Q5").Value = 40

# Row 6
$ws.Range("D6").Value = 44166
$ws.Range("H6").Value = 'Madrigal'
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("N6").Value = '$/caja 40 unidades'
$ws.Range("P6").Value = 250
$ws.Range("Q6").Value = 40

# Row 7
$ws.Range("D7").Value = 44400
$ws.Range("H7").Value = 'Española'
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("P7").Value = 500
$ws.Range("Q7").Value = 30

# Row 8
$ws.Range("D8").Value = 44176
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 11000
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("P8").Value = 275
$ws.Range("Q8").Value = 40

# Row 9
$ws.Range("D9").Value = 44161
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 30
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 275
$ws.Range("Q9").Value = 40

# Row 11
$ws.Range("D11").Value = 44407
$ws.Range("H11").Value = 'Española'
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 18000
$ws.Range("N11").Value = '$/caja 30 unidades'
$ws.Range("P11").Value = 600
$ws.Range("Q11").Value = 30

# Row 12
$ws.Range("D12").Value = 44488
$ws.Range("H12").Value = 'Madrigal'
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("N12").Value = '$/caja 40 unidades'
$ws.Range("P12").Value = 300
$ws.Range("Q12").Value = 40

# Row 13
$ws.Range("D13").Value = 44421
$ws.Range("K13").Value = 16500
$ws.Range("L13").Value = 16500
$ws.Range("M13").Value = 16500
$ws.Range("P13").Value = 550

# Row 14
$ws.Range("D14").Value = 44491
$ws.Range("J14").Value = 200

# Row 15
$ws.Range("D15").Value = 44162
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 10000
$ws.Range("P15").Value = 250

# Row 16
$ws.Range("D16").Value = 44386
$ws.Range("H16").Value = 'Española'
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("N16").Value = '$/caja 30 unidades'
$ws.Range("P16").Value = 500
$ws.Range("Q16").Value = 30

# Row 17
$ws.Range("D17").Value = 44481
$ws.Range("H17").Value = 'Madrigal'
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 11000
$ws.Range("M17").Value = 11000
$ws.Range("N17").Value = '$/caja 50 unidades'
$ws.Range("P17").Value = 220
$ws.Range("Q17").Value = 50
